# Applies the edits described by the target diff:
#  - Swap the Subtype_1 (B) and Mintage (E) values between rows 3 and 4
#  - Swap the Subtype_1 (B) values between rows 10 and 11
#  - Update the active selection on the "2€" sheet to the full row 21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# --- Swap row 3 / row 4 data (columns B and E) ---
$b3 = $ws.Range("B3").Value()
$b4 = $ws.Range("B4").Value()
$ws.Range("B3").Value = $b4
$ws.Range("B4").Value = $b3

$e3 = $ws.Range("E3").Value()
$e4 = $ws.Range("E4").Value()
$ws.Range("E3").Value = $e4
$ws.Range("E4").Value = $e3

# --- Swap row 10 / row 11 data (column B) ---
$b10 = $ws.Range("B10").Value()
$b11 = $ws.Range("B11").Value()
$ws.Range("B10").Value = $b11
$ws.Range("B11").Value = $b10

# --- Update the active cell / selection to the full row 21 ---
$ws.Activate()
$ws.Range("A21:XFD21").Select()
